$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so Excel does not
# auto-coerce numeric-looking strings (e.g. "1.00", "8.02") into numbers,
# matching the source data which stores these as plain text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.419.84'
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").Value = '3.384.58'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '573.97'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").Value = '152.30'
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.483'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = '8.02'
$ws.Range("E9").Value = '  +2.80%  '
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("D12").Value = '3.969.45'
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").Value = '28.47'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '3.382.32'
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").Value = '61.596.67'
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '6.38'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").Value = '14.31'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("D21").Value = '376.25'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").Value = '0.568'
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("D23").Value = '75.58'
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '3.550.01'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").Value = '0.0000110'
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = '7.51'
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("D33").Value = '23.08'
$ws.Range("E33").Value = '  -0.71%  '
$ws.Range("D34").Value = '1.28'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("D35").Value = '5.43'
$ws.Range("E35").Value = '  +2.13%  '
$ws.Range("D36").Value = '169.75'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("D38").Value = '6.85'
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("D39").Value = '30.44'
$ws.Range("E39").Value = '  -3.65%  '
$ws.Range("D40").Value = '3.425.23'
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("E41").Value = '  +1.07%  '
$ws.Range("D42").Value = '42.48'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("E45").Value = '  -5.30%  '
$ws.Range("E46").Value = '  -3.61%  '
$ws.Range("D47").Value = '2.533.44'
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Value = '23.04'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("D49").Value = '6.77'
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '0.0263'
$ws.Range("E51").Value = '  -1.04%  '
